# Daily attendance processing - 2026-01-15 16:46:28
# Reorder the "Recorded By" (column G) contributor list for each session row:
# reverse the order of the comma-separated names/emails recorded in the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $reversed = $parts[($parts.Count - 1)..0]
            $cell.Value = $reversed -join ", "
        }
    }
}
